# Fruta / hortaliza, semanal
# Re-sync the weekly price rows with updated upstream data.
# Applies the per-cell deltas described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44221
$ws.Range("H2").Value = "Cultivar XV región"
$ws.Range("J2").Value = 140
$ws.Range("K2").Value = 5000
$ws.Range("L2").Value = 6000
$ws.Range("M2").Value = 5500
$ws.Range("N2").Value = "`$/caja 10 kilos"
$ws.Range("O2").Value = "Región de Arica y Parinacota"
$ws.Range("P2").Value = 550
$ws.Range("Q2").Value = 10

# Row 5
$ws.Range("D5").Value = 44405
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 140

# Row 6
$ws.Range("D6").Value = 44211
$ws.Range("H6").Value = "Cultivar XV región"
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 140
$ws.Range("K6").Value = 4500
$ws.Range("L6").Value = 5000
$ws.Range("M6").Value = 4750
$ws.Range("N6").Value = "`$/caja 10 kilos"
$ws.Range("O6").Value = "Región de Arica y Parinacota"
$ws.Range("P6").Value = 475
$ws.Range("Q6").Value = 10

# Row 7
$ws.Range("D7").Value = 44412
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 150

# Row 8
$ws.Range("D8").Value = 44454
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 160
$ws.Range("K8").Value = 19000
$ws.Range("L8").Value = 20000
$ws.Range("M8").Value = 19500
$ws.Range("P8").Value = 1083

# Row 10
$ws.Range("D10").Value = 44435
$ws.Range("J10").Value = 100

# Row 11
$ws.Range("D11").Value = 44435
$ws.Range("I11").Value = "Tercera"
$ws.Range("J11").Value = 120
$ws.Range("K11").Value = 14000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 14500
$ws.Range("P11").Value = 806

# Row 12
$ws.Range("D12").Value = 44377
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 17000
$ws.Range("L12").Value = 18000
$ws.Range("M12").Value = 17600
$ws.Range("P12").Value = 978

# Row 13
$ws.Range("D13").Value = 44433
$ws.Range("H13").Value = "Cultivar IV Región"
$ws.Range("I13").Value = "Segunda"
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 17000
$ws.Range("L13").Value = 18000
$ws.Range("M13").Value = 17500
$ws.Range("N13").Value = "`$/bandeja 18 kilos"
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 972
$ws.Range("Q13").Value = 18

# Row 14
$ws.Range("D14").Value = 44433
$ws.Range("I14").Value = "Tercera"
$ws.Range("J14").Value = 120

# Row 15
$ws.Range("D15").Value = 44363
$ws.Range("H15").Value = "Cultivar IV Región"
$ws.Range("I15").Value = "Primera"
$ws.Range("K15").Value = 14000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 14500
$ws.Range("N15").Value = "`$/bandeja 18 kilos"
$ws.Range("O15").Value = "Provincia de Limarí"
$ws.Range("P15").Value = 806
$ws.Range("Q15").Value = 18
